# Add a "Currency" column to the Transactions sheet, between "Amount" (D)
# and "Add to Splitwise? (TRUE/FALSE)" (old E, now F), and select a cell
# on each sheet to match the authored selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")
$cfg = $wb.Worksheets.Item("Config")

# Insert a new column before column E, shifting everything right.
$ws.Range("E1").EntireColumn.Insert()

# Header text + style for the new column (match style of D1 / Amount header).
$ws.Range("E1").Value = "Currency"
$ws.Range("E1").Style = $ws.Range("D1").Style

# Match the column width/format behaviour of the new Currency column (8.85546875, bestFit like column D).
$ws.Columns("E").ColumnWidth = 8.85546875

# Update selections to match the authored state.
$ws.Range("E9").Select()
$cfg.Range("B2:B3").Select()
$cfg.Range("B3").Activate()
